$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "PERIOD TO EXPIRE" (column H) values for rows 3-29, each one day less
# than before (as of the new "LAST UPDATE" date of 04-Nov-2025).
$newPeriod = @{
    3  = 192
    4  = 182
    5  = 190
    6  = 191
    7  = 359
    8  = 204
    9  = 202
    10 = 198
    11 = 189
    12 = 351
    13 = 174
    14 = 300
    15 = 350
    16 = 205
    17 = 310
    18 = 687
    19 = 371
    20 = 709
    21 = 709
    22 = 394
    23 = 36
    24 = -175
    25 = 125
    26 = 128
    27 = 140
    28 = 183
    29 = 612
}

for ($r = 3; $r -le 29; $r++) {
    # Update PERIOD TO EXPIRE (column H) with the new numeric value.
    $ws.Cells($r, 8).Value = $newPeriod[$r]

    # Update LAST UPDATE (column I) to 04-Nov-2025, keeping the cell as
    # plain text (not an auto-converted date serial) and preserving the
    # existing cell style. Writing the value as a text formula first and
    # then pasting-special as values-only avoids Excel's automatic
    # recognition of "dd-mmm-yyyy" look-alike strings as dates.
    $cell = $ws.Cells($r, 9)
    $cell.Formula = "=""04-Nov-2025"""
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$excel.CutCopyMode = $false
